$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note row: record that the log was updated for MATLAB/Simulink R2022a
$ws.Cells.Item(152, 1).Value = "Updated to MATLAB/Simulink R2022a"

# Row 154: Ramps2 trial
$ws.Cells.Item(154, 1).Value = 44641
$ws.Cells.Item(154, 2).Value = "Ramps2"
$ws.Cells.Item(154, 3).Value = 6
$ws.Cells.Item(154, 4).Value = 4
$ws.Cells.Item(154, 5).Value = 4

# Rows 155-160: Sine2 trials
$ws.Cells.Item(155, 1).Value = 44641
$ws.Cells.Item(155, 2).Value = "Sine2"
$ws.Cells.Item(155, 3).Value = 2
$ws.Cells.Item(155, 4).Value = 0.5
$ws.Cells.Item(155, 5).Value = 1
$ws.Cells.Item(155, 6).Value = 0.75

$ws.Cells.Item(156, 1).Value = 44641
$ws.Cells.Item(156, 2).Value = "Sine2"
$ws.Cells.Item(156, 3).Value = 3
$ws.Cells.Item(156, 4).Value = 0.5
$ws.Cells.Item(156, 5).Value = 1
$ws.Cells.Item(156, 6).Value = 1

$ws.Cells.Item(157, 1).Value = 44641
$ws.Cells.Item(157, 2).Value = "Sine2"
$ws.Cells.Item(157, 3).Value = 4
$ws.Cells.Item(157, 4).Value = 0.5
$ws.Cells.Item(157, 5).Value = 1
$ws.Cells.Item(157, 6).Value = 1.5

$ws.Cells.Item(158, 1).Value = 44641
$ws.Cells.Item(158, 2).Value = "Sine2"
$ws.Cells.Item(158, 3).Value = 5
$ws.Cells.Item(158, 4).Value = 0.5
$ws.Cells.Item(158, 5).Value = 1
$ws.Cells.Item(158, 6).Value = 2

$ws.Cells.Item(159, 1).Value = 44641
$ws.Cells.Item(159, 2).Value = "Sine2"
$ws.Cells.Item(159, 3).Value = 6
$ws.Cells.Item(159, 4).Value = 0.5
$ws.Cells.Item(159, 5).Value = 1
$ws.Cells.Item(159, 6).Value = 3

$ws.Cells.Item(160, 1).Value = 44641
$ws.Cells.Item(160, 2).Value = "Sine2"
$ws.Cells.Item(160, 3).Value = 7
$ws.Cells.Item(160, 4).Value = 1.5
$ws.Cells.Item(160, 5).Value = 3
$ws.Cells.Item(160, 6).Value = 0.75

# Update the frozen-pane view / active selection to match the new scroll position
$ws.Application.ActiveWindow.ScrollRow = 142
$ws.Range("B161").Select()
